# Update ownership-rate figures per revised source data (Add files via upload)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 2.8
$ws.Range("D3").Value = 4
$ws.Range("D4").Value = 5.2
$ws.Range("D5").Value = 7.1
$ws.Range("D7").Value = 12.8
$ws.Range("D8").Value = 14.5
$ws.Range("D9").Value = 17.6
$ws.Range("D10").Value = 21.8
$ws.Range("D11").Value = 27.2
$ws.Range("F13").Value = 34.1
$ws.Range("F14").Value = 38.4
$ws.Range("F15").Value = 41.3
$ws.Range("F17").Value = 47.3
$ws.Range("F18").Value = 49
$ws.Range("F19").Value = 51.2
$ws.Range("F20").Value = 53.2
$ws.Range("F22").Value = 56.2
$ws.Range("F23").Value = 56.6
$ws.Range("F24").Value = 58.4
$ws.Range("F25").Value = 60.6
$ws.Range("F26").Value = 60.2
$ws.Range("F27").Value = 62.4
$ws.Range("C29").Value = 63.2
$ws.Range("C30").Value = 63.7
$ws.Range("C31").Value = 64.3
$ws.Range("C33").Value = 63.8
$ws.Range("C34").Value = 66.3
$ws.Range("C35").Value = 67.3
$ws.Range("C36").Value = 67.8
$ws.Range("C37").Value = 67.8
$ws.Range("C38").Value = 68
$ws.Range("C39").Value = 69
$ws.Range("C40").Value = 69.2
$ws.Range("C41").Value = 69.2
$ws.Range("C42").Value = 70.3
$ws.Range("C43").Value = 70.1
$ws.Range("B44").Value = 73.9
$ws.Range("B45").Value = 74
$ws.Range("B47").Value = 74.9
$ws.Range("B49").Value = 75.1
$ws.Range("B51").Value = 75.9
$ws.Range("B52").Value = 76.9
$ws.Range("B53").Value = 76.9
$ws.Range("B54").Value = 77.4
$ws.Range("B55").Value = 76.8
$ws.Range("B56").Value = 76.4
$ws.Range("B57").Value = 78
$ws.Range("B59").Value = 75.9
$ws.Range("B60").Value = 76.4
$ws.Range("B61").Value = 73.6
$ws.Range("B62").Value = 75.7
$ws.Range("G68").Value = 62.2
